$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.894167
$ws.Range("H2").Value = 5.682501
$ws.Range("I2").Value = 0.08880272738170709
$ws.Range("J2").Value = 0.08880272738170709
$ws.Range("M2").Value = 4.559506666666667
$ws.Range("N2").Value = 13.67852
$ws.Range("O2").Value = 0.2891443386304051
$ws.Range("P2").Value = 0.289144338630405
$ws.Range("Q2").Value = 8.636467064280001
$ws.Range("R2").Value = 77.72820357852
$ws.Range("S2").Value = 0.02567680587735986
$ws.Range("T2").Value = 0.02567680587735986
$ws.Range("G3").Value = 1.894167
$ws.Range("H3").Value = 5.682501
$ws.Range("I3").Value = 0.08880272738170709
$ws.Range("J3").Value = 0.08880272738170709
$ws.Range("O3").Value = 0.3132052302231143
$ws.Range("P3").Value = 0.3132052302231142
$ws.Range("Q3").Value = 9.355143067974002
$ws.Range("R3").Value = 84.19628761176601
$ws.Range("S3").Value = 0.02781347867402802
$ws.Range("T3").Value = 0.02781347867402802
$ws.Range("G4").Value = 1.894167
$ws.Range("H4").Value = 5.682501
$ws.Range("I4").Value = 0.08880272738170709
$ws.Range("J4").Value = 0.08880272738170709
$ws.Range("M4").Value = 3.218510666666667
$ws.Range("N4").Value = 9.655531999999999
$ws.Range("O4").Value = 0.2041041292672535
$ws.Range("P4").Value = 0.2041041292672534
$ws.Range("Q4").Value = 6.096396693948
$ws.Range("R4").Value = 54.867570245532
$ws.Range("S4").Value = 0.01812500334880061
$ws.Range("T4").Value = 0.01812500334880061
$ws.Range("G5").Value = 1.894167
$ws.Range("H5").Value = 5.682501
$ws.Range("I5").Value = 0.08880272738170709
$ws.Range("J5").Value = 0.08880272738170709
$ws.Range("M5").Value = 3.052024666666667
$ws.Range("N5").Value = 9.156074
$ws.Range("O5").Value = 0.1935463018792272
$ws.Range("P5").Value = 0.1935463018792272
$ws.Range("Q5").Value = 5.781044406786
$ws.Range("R5").Value = 52.029399661074
$ws.Range("S5").Value = 0.0171874394815186
$ws.Range("T5").Value = 0.0171874394815186
$ws.Range("I6").Value = 0.3983810605389457
$ws.Range("J6").Value = 0.3983810605389456
$ws.Range("M6").Value = 4.559506666666667
$ws.Range("N6").Value = 13.67852
$ws.Range("O6").Value = 0.2891443386304051
$ws.Range("P6").Value = 0.289144338630405
$ws.Range("Q6").Value = 38.74436078509778
$ws.Range("R6").Value = 348.6992470658799
$ws.Range("S6").Value = 0.1151896282724128
$ws.Range("T6").Value = 0.1151896282724128
$ws.Range("I7").Value = 0.3983810605389457
$ws.Range("J7").Value = 0.3983810605389456
$ws.Range("O7").Value = 0.3132052302231143
$ws.Range("P7").Value = 0.3132052302231142
$ws.Range("S7").Value = 0.1247750317826289
$ws.Range("T7").Value = 0.1247750317826289
$ws.Range("I8").Value = 0.3983810605389457
$ws.Range("J8").Value = 0.3983810605389456
$ws.Range("M8").Value = 3.218510666666667
$ws.Range("N8").Value = 9.655531999999999
$ws.Range("O8").Value = 0.2041041292672535
$ws.Range("P8").Value = 0.2041041292672534
$ws.Range("Q8").Value = 27.34926113205644
$ws.Range("R8").Value = 246.143350188508
$ws.Range("S8").Value = 0.0813112194778665
$ws.Range("T8").Value = 0.08131121947786647
$ws.Range("I9").Value = 0.3983810605389457
$ws.Range("J9").Value = 0.3983810605389456
$ws.Range("M9").Value = 3.052024666666667
$ws.Range("N9").Value = 9.156074
$ws.Range("O9").Value = 0.1935463018792272
$ws.Range("P9").Value = 0.1935463018792272
$ws.Range("Q9").Value = 25.93454806741178
$ws.Range("R9").Value = 233.410932606706
$ws.Range("S9").Value = 0.07710518100603747
$ws.Range("T9").Value = 0.07710518100603747
$ws.Range("G10").Value = 5.277913
$ws.Range("H10").Value = 15.833739
$ws.Range("I10").Value = 0.2474402042076373
$ws.Range("J10").Value = 0.2474402042076373
$ws.Range("M10").Value = 4.559506666666667
$ws.Range("N10").Value = 13.67852
$ws.Range("O10").Value = 0.2891443386304051
$ws.Range("P10").Value = 0.289144338630405
$ws.Range("Q10").Value = 24.06467950958666
$ws.Range("R10").Value = 216.58211558628
$ws.Range("S10").Value = 0.07154593419618967
$ws.Range("T10").Value = 0.07154593419618965
$ws.Range("G11").Value = 5.277913
$ws.Range("H11").Value = 15.833739
$ws.Range("I11").Value = 0.2474402042076373
$ws.Range("J11").Value = 0.2474402042076373
$ws.Range("O11").Value = 0.3132052302231143
$ws.Range("P11").Value = 0.3132052302231142
$ws.Range("Q11").Value = 26.067200629786
$ws.Range("R11").Value = 234.604805668074
$ws.Range("S11").Value = 0.07749956612530746
$ws.Range("T11").Value = 0.07749956612530745
$ws.Range("G12").Value = 5.277913
$ws.Range("H12").Value = 15.833739
$ws.Range("I12").Value = 0.2474402042076373
$ws.Range("J12").Value = 0.2474402042076373
$ws.Range("M12").Value = 3.218510666666667
$ws.Range("N12").Value = 9.655531999999999
$ws.Range("O12").Value = 0.2041041292672535
$ws.Range("P12").Value = 0.2041041292672534
$ws.Range("Q12").Value = 16.98701928823867
$ws.Range("R12").Value = 152.883173594148
$ws.Range("S12").Value = 0.05050356742551121
$ws.Range("T12").Value = 0.0505035674255112
$ws.Range("G13").Value = 5.277913
$ws.Range("H13").Value = 15.833739
$ws.Range("I13").Value = 0.2474402042076373
$ws.Range("J13").Value = 0.2474402042076373
$ws.Range("M13").Value = 3.052024666666667
$ws.Range("N13").Value = 9.156074
$ws.Range("O13").Value = 0.1935463018792272
$ws.Range("P13").Value = 0.1935463018792272
$ws.Range("Q13").Value = 16.10832066452067
$ws.Range("R13").Value = 144.974885980686
$ws.Range("S13").Value = 0.04789113646062901
$ws.Range("T13").Value = 0.04789113646062901
$ws.Range("G14").Value = 5.660484666666666
$ws.Range("H14").Value = 16.981454
$ws.Range("I14").Value = 0.2653760078717099
$ws.Range("J14").Value = 0.2653760078717099
$ws.Range("M14").Value = 4.559506666666667
$ws.Range("N14").Value = 13.67852
$ws.Range("O14").Value = 0.2891443386304051
$ws.Range("P14").Value = 0.289144338630405
$ws.Range("Q14").Value = 25.80901757423111
$ws.Range("R14").Value = 232.28115816808
$ws.Range("S14").Value = 0.07673197028444273
$ws.Range("T14").Value = 0.07673197028444272
$ws.Range("G15").Value = 5.660484666666666
$ws.Range("H15").Value = 16.981454
$ws.Range("I15").Value = 0.2653760078717099
$ws.Range("J15").Value = 0.2653760078717099
$ws.Range("O15").Value = 0.3132052302231143
$ws.Range("P15").Value = 0.3132052302231142
$ws.Range("Q15").Value = 27.95669225086267
$ws.Range("R15").Value = 251.610230257764
$ws.Range("S15").Value = 0.08311715364114988
$ws.Range("T15").Value = 0.08311715364114987
$ws.Range("G16").Value = 5.660484666666666
$ws.Range("H16").Value = 16.981454
$ws.Range("I16").Value = 0.2653760078717099
$ws.Range("J16").Value = 0.2653760078717099
$ws.Range("M16").Value = 3.218510666666667
$ws.Range("N16").Value = 9.655531999999999
$ws.Range("O16").Value = 0.2041041292672535
$ws.Range("P16").Value = 0.2041041292672534
$ws.Range("Q16").Value = 18.21833027816978
$ws.Range("R16").Value = 163.964972503528
$ws.Range("S16").Value = 0.05416433901507515
$ws.Range("T16").Value = 0.05416433901507514
$ws.Range("G17").Value = 5.660484666666666
$ws.Range("H17").Value = 16.981454
$ws.Range("I17").Value = 0.2653760078717099
$ws.Range("J17").Value = 0.2653760078717099
$ws.Range("M17").Value = 3.052024666666667
$ws.Range("N17").Value = 9.156074
$ws.Range("O17").Value = 0.1935463018792272
$ws.Range("P17").Value = 0.1935463018792272
$ws.Range("Q17").Value = 17.27593882795511
$ws.Range("R17").Value = 155.483449451596
$ws.Range("S17").Value = 0.05136254493104214
$ws.Range("T17").Value = 0.05136254493104214
